# Workshop_SistemaBancario.pptx - "atualização dos slides do workshop"
#
# On slide 2 ("Material do Workshop"):
#   1. Reposition the content picture placeholder.
#   2. Add a rectangle shape with the project's GitHub URL below it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1. Move the picture placeholder to its new location -------------------
$pic = $s.Shapes.Item(3)
$pic.Left = 436.5428346456693
$pic.Top  = 35.94677165354331

# --- 2. Add the new "Retângulo 2" rectangle with the repo URL --------------
$rect = $s.Shapes.AddShape(1, 40.48188976377953, 489.15, 541.5903937007874, 29.081259842519685)
$rect.Name = "Retângulo 2"
$rect.TextFrame.WordWrap = -1
$rect.TextFrame.AutoSize = 1
$rect.TextFrame.TextRange.Text = "https://github.com/edsonmfeitosa/Workshop_Sistema_Bancario"
